$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the customer name from "Lance" to "Levent"
$ws.Range("C2").Value = "Levent"

# Update invoice number
$ws.Range("B2").Value = 3904

# Update total amount
$ws.Range("E2").Value = 0
